$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in rows 2-6 (re-scaled figures)
# Row 2
$ws.Range("D2").Value = 2309
$ws.Range("E2").Value = 213
$ws.Range("F2").Value = 213
$ws.Range("G2").Value = 210
$ws.Range("H2").Value = 167
$ws.Range("I2").Value = 163
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1107
$ws.Range("L2").Value = 265
$ws.Range("M2").Value = 842
$ws.Range("N2").Value = 823
$ws.Range("O2").Value = 19
$ws.Range("P2").Value = 557
$ws.Range("Q2").Value = 247
$ws.Range("R2").Value = -20
$ws.Range("S2").Value = -203
$ws.Range("T2").Value = 19
$ws.Range("U2").Value = 228
$ws.Range("V2").Value = 96
$ws.Range("W2").Value = 9.210000000000001
$ws.Range("X2").Value = 7.25
$ws.Range("Y2").Value = 21.7
$ws.Range("Z2").Value = 15.13
$ws.Range("AA2").Value = 31.46
$ws.Range("AB2").Value = 49.11
$ws.Range("AC2").Value = 147
$ws.Range("AD2").Value = 11.65
$ws.Range("AE2").Value = 742
$ws.Range("AF2").Value = 2.3
$ws.Range("AG2").Value = 45
$ws.Range("AH2").Value = 2.63
$ws.Range("AI2").Value = 30.58
$ws.Range("AJ2").Value = 111133730

# Row 3
$ws.Range("D3").Value = 2164
$ws.Range("E3").Value = 230
$ws.Range("F3").Value = 230
$ws.Range("G3").Value = 176
$ws.Range("H3").Value = 126
$ws.Range("I3").Value = 126
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1155
$ws.Range("L3").Value = 239
$ws.Range("M3").Value = 916
$ws.Range("N3").Value = 897
$ws.Range("O3").Value = 19
$ws.Range("P3").Value = 557
$ws.Range("Q3").Value = 246
$ws.Range("R3").Value = -50
$ws.Range("S3").Value = -78
$ws.Range("T3").Value = 48
$ws.Range("U3").Value = 198
$ws.Range("V3").Value = 69
$ws.Range("W3").Value = 10.6
$ws.Range("X3").Value = 5.82
$ws.Range("Y3").Value = 14.62
$ws.Range("Z3").Value = 11.13
$ws.Range("AA3").Value = 26.05
$ws.Range("AB3").Value = 62.24
$ws.Range("AC3").Value = 113
$ws.Range("AD3").Value = 12.72
$ws.Range("AE3").Value = 809
$ws.Range("AF3").Value = 1.78
$ws.Range("AG3").Value = 44
$ws.Range("AH3").Value = 3.06
$ws.Range("AI3").Value = 38.79
$ws.Range("AJ3").Value = 111133730

# Row 4
$ws.Range("D4").Value = 2184
$ws.Range("E4").Value = 241
$ws.Range("F4").Value = 241
$ws.Range("G4").Value = 233
$ws.Range("H4").Value = 180
$ws.Range("I4").Value = 179
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1245
$ws.Range("L4").Value = 204
$ws.Range("M4").Value = 1041
$ws.Range("N4").Value = 1023
$ws.Range("O4").Value = 18
$ws.Range("P4").Value = 557
$ws.Range("Q4").Value = 117
$ws.Range("R4").Value = 39
$ws.Range("S4").Value = -87
$ws.Range("T4").Value = 22
$ws.Range("U4").Value = 95
$ws.Range("V4").Value = 36
$ws.Range("W4").Value = 11.04
$ws.Range("X4").Value = 8.25
$ws.Range("Y4").Value = 18.66
$ws.Range("Z4").Value = 15.02
$ws.Range("AA4").Value = 19.6
$ws.Range("AB4").Value = 84.91
$ws.Range("AC4").Value = 161
$ws.Range("AD4").Value = 9.77
$ws.Range("AE4").Value = 922
$ws.Range("AF4").Value = 1.71
$ws.Range("AG4").Value = 55
$ws.Range("AH4").Value = 3.49
$ws.Range("AI4").Value = 34.04
$ws.Range("AJ4").Value = 111133730

# Row 5
$ws.Range("D5").Value = 1881
$ws.Range("E5").Value = 142
$ws.Range("F5").Value = 142
$ws.Range("G5").Value = 126
$ws.Range("H5").Value = 99
$ws.Range("I5").Value = 96
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1253
$ws.Range("L5").Value = 176
$ws.Range("M5").Value = 1077
$ws.Range("N5").Value = 1056
$ws.Range("O5").Value = 21
$ws.Range("P5").Value = 557
$ws.Range("Q5").Value = 142
$ws.Range("R5").Value = -13
$ws.Range("S5").Value = -61
$ws.Range("T5").Value = 28
$ws.Range("U5").Value = 114
$ws.Range("V5").Value = 36
$ws.Range("W5").Value = 7.54
$ws.Range("X5").Value = 5.25
$ws.Range("Y5").Value = 9.210000000000001
$ws.Range("Z5").Value = 7.91
$ws.Range("AA5").Value = 16.36
$ws.Range("AB5").Value = 90.83
$ws.Range("AC5").Value = 86
$ws.Range("AD5").Value = 16.32
$ws.Range("AE5").Value = 952
$ws.Range("AF5").Value = 1.48
$ws.Range("AG5").Value = 43
$ws.Range("AH5").Value = 3.06
$ws.Range("AI5").Value = 49.85
$ws.Range("AJ5").Value = 111133730

# Row 6
$ws.Range("D6").Value = 1836
$ws.Range("E6").Value = 21
$ws.Range("F6").Value = 21
$ws.Range("G6").Value = 27
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = 23
$ws.Range("K6").Value = 1224
$ws.Range("L6").Value = 170
$ws.Range("M6").Value = 1054
$ws.Range("N6").Value = 1035
$ws.Range("P6").Value = 557
$ws.Range("Q6").Value = 52
$ws.Range("R6").Value = -40
$ws.Range("S6").Value = -49
$ws.Range("T6").Value = 15
$ws.Range("U6").Value = 37
$ws.Range("V6").Value = 35
$ws.Range("W6").Value = 1.14
$ws.Range("X6").Value = 1.11
$ws.Range("Y6").Value = 2.21
$ws.Range("Z6").Value = 1.64
$ws.Range("AA6").Value = 16.17
$ws.Range("AB6").Value = 87.17
$ws.Range("AC6").Value = 21
$ws.Range("AD6").Value = 53.01
$ws.Range("AE6").Value = 934
$ws.Range("AF6").Value = 1.18
$ws.Range("AG6").Value = 23
$ws.Range("AH6").Value = 2.09
$ws.Range("AI6").Value = 110.53
$ws.Range("AJ6").Value = 111133730

# Clear numeric data in rows 7-9 (projection rows removed)
# Row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

